$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 1383.0965220845178
$ws.Range("C1").Value = 1428.823605919041
$ws.Range("B2").Value = 1482.0442259078409
$ws.Range("C2").Value = 1329.6420705791181
$ws.Range("A3").Value = 2343.9137069531043
$ws.Range("B3").Value = 1585.8755819885312
$ws.Range("C3").Value = 1454.0869954572026
$ws.Range("A4").Value = 2320.9625000851147
$ws.Range("B4").Value = 1780.0461948646082
$ws.Range("C4").Value = 1731.5928878222501
$ws.Range("A5").Value = 2422.6595198945024
$ws.Range("B5").Value = 1674.9654899696432
$ws.Range("C5").Value = 1622.5281738093465
$ws.Range("A6").Value = 2360.7084402304272
$ws.Range("B6").Value = 1774.6285757658027
$ws.Range("C6").Value = 1784.4329783711232
$ws.Range("A7").Value = 1992.9848472506942
$ws.Range("B7").Value = 1566.2947458768306
$ws.Range("C7").Value = 1482.2946237954445
$ws.Range("A8").Value = 2135.7478628450635
$ws.Range("B8").Value = 1655.7390016109439
$ws.Range("C8").Value = 1507.3584903465025
$ws.Range("A9").Value = 2471.0184236373766
$ws.Range("B9").Value = 1806.1208547647914
$ws.Range("C9").Value = 1570.5840636497583
$ws.Range("A10").Value = 2111.5944043737513
$ws.Range("B10").Value = 1384.8197897029249
$ws.Range("C10").Value = 1313.1822179032408
$ws.Range("A11").Value = 1970.0049584060432
$ws.Range("B11").Value = 1414.6269501252143
$ws.Range("C11").Value = 1312.1828705205139
$ws.Range("A12").Value = 2787.9537156936608
$ws.Range("B12").Value = 2270.4398262776958
$ws.Range("C12").Value = 2036.7362819098314
$ws.Range("A13").Value = 2315.1858317581732
$ws.Range("B13").Value = 1777.4423331486139
$ws.Range("C13").Value = 1794.9018703913953
$ws.Range("A14").Value = 2593.0482642470733
$ws.Range("B14").Value = 1924.4781793533441
$ws.Range("C14").Value = 1698.5718122336648
$ws.Range("A15").Value = 2508.5108803272624
$ws.Range("B15").Value = 2026.7228967099065
$ws.Range("C15").Value = 1827.5083266369336
$ws.Range("A16").Value = 2205.5211166519052
$ws.Range("B16").Value = 1546.5982352544029
$ws.Range("C16").Value = 1276.461477199553
$ws.Range("A17").Value = 2225.5448130781733
$ws.Range("B17").Value = 1682.2603177883561
$ws.Range("C17").Value = 1573.6650186867739
$ws.Range("A18").Value = 2487.9072724582679
$ws.Range("B18").Value = 2063.1343444789436
$ws.Range("C18").Value = 1909.7693489385661
$ws.Range("A19").Value = 1761.81626294991
$ws.Range("B19").Value = 1928.7684870713738
$ws.Range("C19").Value = 1884.2930518248493
$ws.Range("A20").Value = 2351.2568365245052
$ws.Range("B20").Value = 1850.0736691983486
$ws.Range("C20").Value = 1653.3462609483511
$ws.Range("A21").Value = 2583.1100843929667
$ws.Range("B21").Value = 1907.0148615932583
$ws.Range("C21").Value = 1813.087800925716
$ws.Range("A22").Value = 2448.3691807967398
$ws.Range("B22").Value = 1890.2457180811034
$ws.Range("C22").Value = 1644.2810920719
